$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.906.85'
$ws.Range("E2").Value = '  +3.22%  '
$ws.Range("D3").Value = '3.036.47'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.49'
$ws.Range("E5").Value = '  +1.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.23'
$ws.Range("E6").Value = '  +6.81%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.032.02'
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("E9").Value = '  +0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.35'
$ws.Range("E10").Value = '  +10.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.151'
$ws.Range("E11").Value = '  +6.33%  '
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'
$ws.Range("E13").Value = '  +4.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.93'
$ws.Range("E14").Value = '  +3.24%  '
$ws.Range("D16").Value = '3.541.61'
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").Value = '62.922.00'
$ws.Range("E17").Value = '  +3.14%  '
$ws.Range("E18").Value = '  +0.69%  '
$ws.Range("D19").Value = '3.038.89'
$ws.Range("E19").Value = '  +2.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '455.84'
$ws.Range("E20").Value = '  +2.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.22'
$ws.Range("E21").Value = '  +2.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.692'
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.49'
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.99'
$ws.Range("E24").Value = '  +2.24%  '
$ws.Range("E25").Value = '  +5.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.82'
$ws.Range("E26").Value = '  +9.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.19'
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  +2.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.41'
$ws.Range("E30").Value = '  +8.85%  '
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  +5.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.65'
$ws.Range("E33").Value = '  +2.11%  '
$ws.Range("E34").Value = '  +4.61%  '
$ws.Range("D35").Value = '0.0₃0855'
$ws.Range("E35").Value = '  +10.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.03'
$ws.Range("E36").Value = '  +2.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.91'
$ws.Range("E37").Value = '  +3.53%  '
$ws.Range("E38").Value = '  +14.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.10'
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.49'
$ws.Range("E40").Value = '  +0.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.10'
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E42").Value = '  +5.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.295'
$ws.Range("E43").Value = '  +12.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.19'
$ws.Range("E44").Value = '  +10.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '392.88'
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0356'
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("D47").Value = '2.746.08'
$ws.Range("E47").Value = '  +2.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.77'
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.20'
$ws.Range("E50").Value = '  +3.04%  '
$ws.Range("E51").Value = '  +0.97%  '
